# Weekly price update: a new week's worth of price rows (2 rows) was inserted
# at the top of the "Acelga" price history block (which starts at row 534),
# pushing all subsequent rows down by two. The two brand-new rows (534-535)
# duplicate the field values of the rows that end up right below them
# (536-537, which are what used to be rows 534-535) except for the handful
# of cells that actually hold new data (Fecha / Volumen / Precio minimo /
# Precio promedio ponderado / Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 534; everything from the old row 534
# downward (including the two brand new rows that will land at the very
# bottom, 641-642) shifts down by two rows. Excel carries the row-534
# formatting (incl. the date style on column D) onto the freshly inserted
# rows automatically.
$ws.Rows("534:535").Insert()

# The content that used to live in rows 534-535 is now in rows 536-537;
# clone it into the two new rows so every column (A,B,C,E,F,G,H,I,L,N,O,Q,R)
# starts out identical, then we only need to touch the handful of cells that
# actually differ per the new data.
$ws.Range("A536:R536").Copy()
$ws.Range("A534").PasteSpecial()
$ws.Range("A537:R537").Copy()
$ws.Range("A535").PasteSpecial()
$excel.CutCopyMode = 0

# Row 534 (was a clone of the old row 534 / now row 536): update the cells
# that hold genuinely new data.
$ws.Range("D534").Value = 45258
$ws.Range("J534").Value = 300
$ws.Range("K534").Value = 700
$ws.Range("M534").Value = 700
$ws.Range("P534").Value = 700

# Row 535 (was a clone of the old row 535 / now row 537): update the cells
# that hold genuinely new data.
$ws.Range("D535").Value = 45258
$ws.Range("J535").Value = 200
